$wb = $excel.ActiveWorkbook

# --- Select A1:B1 on the previously-last sheet (AddCattletoUserdefinedShed) ---
# Do this BEFORE creating/activating the new sheet so that sheet loses
# tabSelected once the new sheet becomes active.
$prevLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$prevLast.Range("A1:B1").Select()

# --- Add the new sheet after the current last sheet ---
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prevLast)
$ws.Name = "DeleteUserdefinedGroup"

# --- Fill in cell values in the order that reproduces the shared-string
#     append order from the original authoring session ---
$ws.Range("B2").Value = "Are you sure you want to delete the Group grp_No from Shed shed_No?"
$ws.Range("B3").Value = "Group grp_No is the last group in the Shed shed_No. Deleting this Group grp_No will also delete the Shed shed_No. Do you want to proceed?"
$ws.Range("C1").Value = "successMessage"
$ws.Range("A2").Value = "DeleteUserdefinedGroup_deleteGrp"
$ws.Range("A3").Value = "DeleteUserdefinedGroup_deletelastGrp"
$ws.Range("A4").Value = "DeleteUserdefinedGroup_deleteGrpwithCattle"
$ws.Range("B4").Value = "Move all cattle out from this group"
$ws.Range("C2").Value = "Group - grp_No has been deleted from Shed shed_No."
$ws.Range("B5").Value = "Group grp_No is the last group in the last active Shed shed_No. Deleting this Group grp_No will delete all Sheds and Groups in Housing and will Reset Housing. Do you want to proceed?"
$ws.Range("A5").Value = "DeleteUserdefinedGroup_deletelastGrpoflastShed"
$ws.Range("C5").Value = "All cattle been moved out from Group grp_No"
$ws.Range("C3").Value = "Group - grp_No has been deleted from Shed shed_No."
$ws.Range("A1").Value = "testcase"
$ws.Range("B1").Value = "warningMessage"

# --- Formatting: reuse the bold+border header style (index 1, same as the
#     other sheets) by copying format from an existing styled cell, then add
#     wrap-text on B1 to get the combined bold+border+wrap style ---
$headerStyleSrc = $wb.Worksheets.Item(13).Range("A1")
$headerStyleSrc.Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").WrapText = $true

$ws.Range("B2:B5").WrapText = $true

# --- Column widths / row heights ---
$ws.Columns("A").ColumnWidth = 46.5
$ws.Columns("B").ColumnWidth = 83
$ws.Columns("C").ColumnWidth = 52.1

$ws.Rows(3).RowHeight = 30
$ws.Rows(5).RowHeight = 28.5

# --- Selection / active cell for the new sheet ---
$ws.Range("B11").Select()

Write-Output "done"
